# "Refined metadata to be additional tab"
#
# 1) Update the "time_taken" timestamps (column F) on the existing "data"
#    sheet to reflect the re-run of the panel query.
# 2) Add a new "metadata" worksheet (placed after "data") describing the
#    panel query itself (name/id/version/retrieval time/request url).

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)
$data.Name = "data"

# --- 1) refresh the per-row query timestamps on the "data" sheet -----------
$timestamps = @(
    "2021-10-05 14:33:51.097477",
    "2021-10-05 14:33:51.097484",
    "2021-10-05 14:33:51.097487",
    "2021-10-05 14:33:51.097490",
    "2021-10-05 14:33:51.097493",
    "2021-10-05 14:33:51.097495",
    "2021-10-05 14:33:51.097498",
    "2021-10-05 14:33:51.097500",
    "2021-10-05 14:33:51.097503",
    "2021-10-05 14:33:51.097505",
    "2021-10-05 14:33:51.097508",
    "2021-10-05 14:33:51.097523",
    "2021-10-05 14:33:51.097527",
    "2021-10-05 14:33:51.097530",
    "2021-10-05 14:33:51.097532",
    "2021-10-05 14:33:51.097534",
    "2021-10-05 14:33:51.097537",
    "2021-10-05 14:33:51.097539",
    "2021-10-05 14:33:51.097541",
    "2021-10-05 14:33:51.097544",
    "2021-10-05 14:33:51.097546",
    "2021-10-05 14:33:51.097549",
    "2021-10-05 14:33:51.097551",
    "2021-10-05 14:33:51.097553",
    "2021-10-05 14:33:51.097556",
    "2021-10-05 14:33:51.097558",
    "2021-10-05 14:33:51.097561",
    "2021-10-05 14:33:51.097563",
    "2021-10-05 14:33:51.097565",
    "2021-10-05 14:33:51.097568",
    "2021-10-05 14:33:51.097570",
    "2021-10-05 14:33:51.097572"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2) add the "metadata" tab, placed after "data" -------------------------
$metadata = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$metadata.Name = "metadata"

# Reuse the "data" sheet's existing header / first-column formatting
# (bold + thin border + centered-top alignment) instead of re-building
# new style entries cell-by-cell.
$data.Range("B1:F1").Copy()
$metadata.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$metadata.Range("G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$metadata.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $metadata.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$metadata.Cells.Item(2, 1).Value = 0
$metadata.Cells.Item(2, 2).Value = "Fatty Acid Oxidation Defects"
$metadata.Cells.Item(2, 3).Value = 103
# data_version ("1.1") must stay text, not be coerced to a float
$metadata.Cells.Item(2, 4).NumberFormat = "@"
$metadata.Cells.Item(2, 4).Value = "1.1"
$metadata.Cells.Item(2, 5).Value = "2021-03-04T02:55:50.885198Z"
$metadata.Cells.Item(2, 6).Value = "2021-10-05 14:33:51.094341"
$metadata.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/103/?format=json"

$data.Select()

Write-Output "done"
